$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A (row labels) so the longer labels are fully visible
$ws.Columns.Item(1).ColumnWidth = 19.21875

# Corrected normalization: Z (Total) should only sum the "idba"/"spades"
# assembler totals in columns C, P, U and Y instead of summing every
# category column C:Y.
$ws.Range("Z2").Formula = "=C2+P2+U2+Y2"
$ws.Range("Z3:Z19").Formula = "=C3+P3+U3+Y3"

# Restore the active selection to U1
$ws.Range("U1").Select()
